$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.204.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.23%  '
$ws.Range("D3").Value = "'2.759.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.06%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'579.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = "'154.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.57%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  +1.36%  '
$ws.Range("D9").Value = "'2.758.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.36%  '
$ws.Range("E10").Value = '  +1.97%  '
$ws.Range("E11").Value = '  +4.24%  '
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("E13").Value = '  +3.41%  '
$ws.Range("D14").Value = "'3.246.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.93%  '
$ws.Range("D15").Value = "'26.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.86%  '
$ws.Range("D16").Value = "'64.106.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.09%  '
$ws.Range("E17").Value = '  +6.02%  '
$ws.Range("D18").Value = "'2.755.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.31%  '
$ws.Range("D19").Value = "'11.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.42%  '
$ws.Range("E20").Value = '  +2.41%  '
$ws.Range("D21").Value = "'360.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.85%  '
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  -1.36%  '
$ws.Range("D25").Value = "'66.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.95%  '
$ws.Range("E26").Value = '  +4.94%  '
$ws.Range("E27").Value = '  +4.86%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").Value = "'0.0₃0911"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.54%  '
$ws.Range("D30").Value = "'2.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").Value = "'7.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.41%  '
$ws.Range("D32").Value = "'1.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +17.04%  '
$ws.Range("D33").Value = "'172.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.26%  '
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").Value = "'20.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.49%  '
$ws.Range("D36").Value = "'4.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.52%  '
$ws.Range("E37").Value = '  +8.06%  '
$ws.Range("E38").Value = '  +9.24%  '
$ws.Range("D39").Value = "'1.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +15.09%  '
$ws.Range("D40").Value = "'347.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.87%  '
$ws.Range("E41").Value = '  +5.22%  '
$ws.Range("D42").Value = "'39.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("E43").Value = '  +8.88%  '
$ws.Range("D44").Value = "'21.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.42%  '
$ws.Range("D45").Value = "'21.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.97%  '
$ws.Range("E46").Value = '  +4.27%  '
$ws.Range("D47").Value = "'0.647"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.07%  '
$ws.Range("D48").Value = "'137.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("E49").Value = '  +2.35%  '
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("E51").Value = '  +0.24%  '